$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @{Row=2; A=45919.01041666666; B=1180},
    @{Row=3; A=45919.02083333334; B=1179},
    @{Row=4; A=45919.03125; B=1183},
    @{Row=5; A=45919.04166666666; B=1176},
    @{Row=6; A=45919.05208333334; B=1232},
    @{Row=7; A=45919.0625; B=1236},
    @{Row=8; A=45919.07291666666; B=1244},
    @{Row=9; A=45919.08333333334; B=1244},
    @{Row=10; A=45919.09375; B=1237},
    @{Row=11; A=45919.10416666666; B=1229},
    @{Row=12; A=45919.11458333334; B=1234},
    @{Row=13; A=45919.125; B=1238},
    @{Row=14; A=45919.13541666666; B=1191},
    @{Row=15; A=45919.14583333334; B=1187},
    @{Row=16; A=45919.15625; B=1181},
    @{Row=17; A=45919.16666666666; B=1176},
    @{Row=18; A=45919.17708333334; B=1064},
    @{Row=19; A=45919.1875; B=1061},
    @{Row=20; A=45919.19791666666; B=1058},
    @{Row=21; A=45919.20833333334; B=1057},
    @{Row=22; A=45919.21875; B=882},
    @{Row=23; A=45919.22916666666; B=883},
    @{Row=24; A=45919.23958333334; B=885},
    @{Row=25; A=45919.25; B=886},
    @{Row=26; A=45919.26041666666; B=758},
    @{Row=27; A=45919.27083333334; B=763},
    @{Row=28; A=45919.28125; B=761},
    @{Row=29; A=45919.29166666666; B=759},
    @{Row=30; A=45919.30208333334; B=682},
    @{Row=31; A=45919.3125; B=682},
    @{Row=32; A=45919.32291666666; B=683},
    @{Row=33; A=45919.33333333334; B=681},
    @{Row=34; A=45919.34375; B=560},
    @{Row=35; A=45919.35416666666; B=561},
    @{Row=36; A=45919.36458333334; B=563},
    @{Row=37; A=45919.375; B=562},
    @{Row=38; A=45919.38541666666; B=537},
    @{Row=39; A=45919.39583333334; B=541},
    @{Row=40; A=45919.40625; B=545},
    @{Row=41; A=45919.41666666666; B=549},
    @{Row=42; A=45919.42708333334; B=639},
    @{Row=43; A=45919.4375; B=642},
    @{Row=44; A=45919.44791666666; B=645},
    @{Row=45; A=45919.45833333334; B=647},
    @{Row=46; A=45919.46875; B=730},
    @{Row=47; A=45919.47916666666; B=730},
    @{Row=48; A=45919.48958333334; B=731},
    @{Row=49; A=45919.5; B=732},
    @{Row=50; A=45919.51041666666; B=717},
    @{Row=51; A=45919.52083333334; B=717},
    @{Row=52; A=45919.53125; B=716},
    @{Row=53; A=45919.54166666666; B=715},
    @{Row=54; A=45919.55208333334; B=700},
    @{Row=55; A=45919.5625; B=701},
    @{Row=56; A=45919.57291666666; B=701},
    @{Row=57; A=45919.58333333334; B=701},
    @{Row=58; A=45919.59375; B=785},
    @{Row=59; A=45919.60416666666; B=785},
    @{Row=60; A=45919.61458333334; B=784},
    @{Row=61; A=45919.625; B=784},
    @{Row=62; A=45919.63541666666; B=767},
    @{Row=63; A=45919.64583333334; B=765},
    @{Row=64; A=45919.65625; B=764},
    @{Row=65; A=45919.66666666666; B=763},
    @{Row=66; A=45919.67708333334; B=671},
    @{Row=67; A=45919.6875; B=666},
    @{Row=68; A=45919.69791666666; B=665},
    @{Row=69; A=45919.70833333334; B=660},
    @{Row=70; A=45919.71875; B=480},
    @{Row=71; A=45919.72916666666; B=474},
    @{Row=72; A=45919.73958333334; B=463},
    @{Row=73; A=45919.75; B=457},
    @{Row=74; A=45919.76041666666; B=298},
    @{Row=75; A=45919.77083333334; B=291},
    @{Row=76; A=45919.78125; B=290},
    @{Row=77; A=45919.79166666666; B=289},
    @{Row=78; A=45919.80208333334; B=260},
    @{Row=79; A=45919.8125; B=259},
    @{Row=80; A=45919.82291666666; B=259},
    @{Row=81; A=45919.83333333334; B=258},
    @{Row=82; A=45919.84375; B=267},
    @{Row=83; A=45919.85416666666; B=267},
    @{Row=84; A=45919.86458333334; B=267},
    @{Row=85; A=45919.875; B=267},
    @{Row=86; A=45919.88541666666; B=314},
    @{Row=87; A=45919.89583333334; B=315},
    @{Row=88; A=45919.90625; B=315},
    @{Row=89; A=45919.91666666666; B=315},
    @{Row=90; A=45919.92708333334; B=375},
    @{Row=91; A=45919.9375; B=376},
    @{Row=92; A=45919.94791666666; B=377},
    @{Row=93; A=45919.95833333334; B=378},
    @{Row=94; A=45919.96875; B=$null},
    @{Row=95; A=45919.97916666666; B=$null},
    @{Row=96; A=45919.98958333334; B=$null},
    @{Row=97; A=45920; B=$null}
)

foreach ($item in $rowsData) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    if ($null -ne $item.B) {
        $ws.Cells.Item($r, 2).Value = $item.B
    }
}
